$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write new translation rows (19-33) for the CONTACT.* keys ---
# The write order below matches the first-appearance order of each new
# shared string in the target workbook, so the regenerated sharedStrings.xml
# table lines up with the authored edit (Excel de-dupes by first use).
$ws.Range('A19').Value = 'CONTACT.INPUT_EMPTY'
$ws.Range('A20').Value = 'CONTACT.EMAIL_UNVALID'
$ws.Range('B20').Value = 'Email non valide'
$ws.Range('C20').Value = 'Email unvalid'
$ws.Range('A21').Value = 'CONTACT.TITLE'
$ws.Range('A22').Value = 'CONTACT.SUPPORT_MESSAGE'
$ws.Range('A24').Value = 'CONTACT.ADMIN_MESSAGE'
$ws.Range('A23').Value = 'CONTACT.SUPPORT_TITLE'
$ws.Range('A25').Value = 'CONTACT.ADMIN_TITLE'
$ws.Range('B22').Value = 'Vous rencontrez un problème technique lié à l''application.'
$ws.Range('C22').Value = 'You have a technical problem with the application.'
$ws.Range('B23').Value = 'Contacter le support'
$ws.Range('C23').Value = 'Contact support'
$ws.Range('C24').Value = 'You have a functional problem.'
$ws.Range('B24').Value = 'Vous rencontrez un problème fonctionnel.'
$ws.Range('B25').Value = 'Contacter un administrateur'
$ws.Range('C25').Value = 'Contact an administrator'
$ws.Range('A26').Value = 'CONTACT.BACK_BUTTON'
$ws.Range('A27').Value = 'CONTACT.INPUT_LOGIN'
$ws.Range('A28').Value = 'CONTACT.INPUT_NAME'
$ws.Range('A29').Value = 'CONTACT.INPUT_EMAIL'
$ws.Range('A30').Value = 'CONTACT.INPUT_FIRSTNAME'
$ws.Range('A31').Value = 'CONTACT.INPUT_PHONE'
$ws.Range('A32').Value = 'CONTACT.INPUT_MESSAGE'
$ws.Range('A33').Value = 'CONTACT.INPUT_SUBMIT'
$ws.Range('B26').Value = 'Retour'
$ws.Range('C26').Value = 'Back'
$ws.Range('C28').Value = 'Name'
$ws.Range('B29').Value = 'Email'
$ws.Range('C30').Value = 'Firstname'
$ws.Range('B28').Value = 'Nom'
$ws.Range('B30').Value = 'Prénom'
$ws.Range('B31').Value = 'Téléphone (optionnel)'
$ws.Range('C31').Value = 'Phone (optionnal)'
$ws.Range('B33').Value = 'Envoyer'
$ws.Range('C33').Value = 'Send'
$ws.Range('B32').Value = 'Message (300 caractères)'
$ws.Range('C32').Value = 'Message (300 characters)'

# --- Remaining cells that duplicate an already-created shared string ---
$ws.Range('B19').Value = 'Champs requis'
$ws.Range('C19').Value = 'Field required'
$ws.Range('B21').Value = 'Contact / Aide'
$ws.Range('C21').Value = 'Contact / Help'
$ws.Range('B27').Value = 'Login'
$ws.Range('C27').Value = 'Login'
$ws.Range('C29').Value = 'Email'

# --- Column widths (A wider for keys, B wider for French text) ---
$ws.Columns.Item(1).ColumnWidth = 30.43
$ws.Columns.Item(2).ColumnWidth = 48.6

# --- Selection matches the author's final cursor position ---
$ws.Range("B30").Select()

Write-Output "done"
